# Insert a new price-report row above the current row 77.
# This pushes the existing rows 77..163 down to 78..164 (Excel's normal
# Insert behaviour), which is exactly the "rotation" seen in the target
# diff: every existing row's data slides down by one row, and the row
# that used to be the very last one (163) ends up duplicated as the new
# last row (164). Row 77 itself is then populated with the brand-new
# record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(77).Insert()

$ws.Range("A77").Value = 5
$ws.Range("B77").Value = "Macroferia Regional de Talca"
$ws.Range("C77").Value = "Maule"
$ws.Range("D77").Value = 44994
$ws.Range("E77").Value = 7
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108002
$ws.Range("J77").Value = "Mango"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 200
$ws.Range("N77").Value = 8000
$ws.Range("O77").Value = 8000
$ws.Range("P77").Value = 8000
$ws.Range("Q77").Value = "$/bandeja 4 kilos"
$ws.Range("R77").Value = "Perú"
$ws.Range("S77").Value = 2000
$ws.Range("T77").Value = 4
